$wb = $excel.ActiveWorkbook

# --- Companies sheet: selection moved to E26, no longer the tab-selected sheet ---
$wb.Worksheets.Item("Companies").Range("E26").Select()

# --- Deals sheet: selection changed to the whole first row ---
$wb.Worksheets.Item("Deals").Range("A1:XFD1").Select()

# --- Add the new "Cases" sheet with its data ---
$ws = $wb.Worksheets.Add()
$ws.Name = "Cases"

# Header row (row 1)
$ws.Cells.Item(1,1).Value = "title"
$ws.Cells.Item(1,2).Value = "status"
$ws.Cells.Item(1,3).Value = "identifier"
$ws.Cells.Item(1,4).Value = "type"
$ws.Cells.Item(1,5).Value = "priority"
$ws.Cells.Item(1,6).Value = "contact"

# Data filled column by column (matches the source-file authoring order)
$ws.Cells.Item(2,1).Value = "CaseTitle1"
$ws.Cells.Item(3,1).Value = "CaseTitle2"

$ws.Cells.Item(2,2).Value = "Awaiting input"
$ws.Cells.Item(3,2).Value = "Enquiring"

$ws.Cells.Item(2,3).Value = "aaaa"
$ws.Cells.Item(3,3).Value = "bbbb"

$ws.Cells.Item(2,4).Value = "Business Support"
$ws.Cells.Item(3,4).Value = "Complaint"

$ws.Cells.Item(2,5).Value = "High"
$ws.Cells.Item(3,5).Value = "Low"

$ws.Cells.Item(2,6).Value = "aaaaa"
$ws.Cells.Item(3,6).Value = "zzzx"

# Header row highlight (yellow fill, same style as the other sheets' header rows)
$ws.Range("A1:F1").Interior.Color = 65535

# Column widths (auto-fit to content)
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()

# Move the new sheet to the end of the tab strip (after "Deals")
$ws.Move()

# Re-fetch by name (post-move) so the selection/activation targets the right sheet
$cases = $wb.Worksheets.Item("Cases")
$cases.Range("D5").Select()
$cases.Activate()
